$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Updated portfolio weights (column N, rows 2-21)
$ws.Range("N2").Value = 0.03391593108946795
$ws.Range("N3").Value = 0.001
$ws.Range("N4").Value = 0.001
$ws.Range("N5").Value = 0.01583706974633639
$ws.Range("N6").Value = 0.15
$ws.Range("N7").Value = 0.001
$ws.Range("N8").Value = 0.06185955134729807
$ws.Range("N9").Value = 0.005090972792413145
$ws.Range("N10").Value = 0.001000000000000019
$ws.Range("N11").Value = 0.15
$ws.Range("N12").Value = 0.02867792884142787
$ws.Range("N13").Value = 0.095365390944626
$ws.Range("N14").Value = 0.001253155238430658
$ws.Range("N16").Value = 0.15
$ws.Range("N17").Value = 0.15
$ws.Range("N19").Value = 0.15
$ws.Range("N20").Value = 0.001000000000000008
$ws.Range("N21").Value = 0.001000000000000014

# Recomputed "Portfolio return ln" row (row 22)
$ws.Range("B22").Value = -0.000915621970130807
$ws.Range("C22").Value = 0.0319937614770624
$ws.Range("D22").Value = 0.009365989214743591
$ws.Range("E22").Value = 0.0177354452815994
$ws.Range("F22").Value = 0.04298969661505471
$ws.Range("G22").Value = 0.06497452557344659
$ws.Range("H22").Value = 0.01501434834603236
$ws.Range("I22").Value = 0.01027818893652233
$ws.Range("J22").Value = 0.007001422568540781
$ws.Range("K22").Value = 0.03473767821508985
$ws.Range("L22").Value = -0.004574972477059985
$ws.Range("M22").Value = 0.02696646619597313

# Recomputed "Portfolio return" row (row 23)
$ws.Range("B23").Value = 0.9990847970837572
$ws.Range("C23").Value = 1.03251106394107
$ws.Range("D23").Value = 1.00940998734646
$ws.Range("E23").Value = 1.017893652197353
$ws.Range("F23").Value = 1.043927138812035
$ws.Range("G23").Value = 1.067131839466267
$ws.Range("H23").Value = 1.015127629913726
$ws.Range("I23").Value = 1.010331190953037
$ws.Range("J23").Value = 1.007025989829323
$ws.Range("K23").Value = 1.035348078817911
$ws.Range("L23").Value = 0.9954354767684461
$ws.Range("M23").Value = 1.027333351790567
$ws.Range("N23").Value = 1.136189565080878
